# Updated symbol list on Tue Jan 24 14:25:29 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking rows. These columns are stored as plain text (e.g.
# "314.36", "3.29%"), so each value is written with a leading apostrophe
# to force text entry and keep Excel from auto-converting the
# number-looking / percent-looking strings into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.36"
$ws.Range("E2").Value = "'3.29%"
$ws.Range("D3").Value = "'36.02"
$ws.Range("E3").Value = "'1.34%"
$ws.Range("D4").Value = "'5.115"
$ws.Range("E4").Value = "'1.15%"
$ws.Range("D5").Value = "'0.08135"
$ws.Range("E5").Value = "'3.39%"
$ws.Range("D6").Value = "'2.132"
$ws.Range("E6").Value = "'1.14%"
$ws.Range("D7").Value = "'4.137"
$ws.Range("E7").Value = "'0.13%"
$ws.Range("D8").Value = "'7.947"
$ws.Range("E8").Value = "'0.29%"
$ws.Range("D9").Value = "'0.9317"
$ws.Range("E9").Value = "'1.12%"
$ws.Range("D10").Value = "'0.1040"
$ws.Range("E10").Value = "'6.36%"
$ws.Range("E11").Value = "'5.25%"
$ws.Range("D12").Value = "'0.09147"
$ws.Range("E12").Value = "'6.15%"
$ws.Range("D13").Value = "'0.03634"
$ws.Range("E13").Value = "'2.81%"
$ws.Range("D14").Value = "'0.09895"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("D15").Value = "'0.001431"
$ws.Range("E15").Value = "'-0.10%"
$ws.Range("D16").Value = "'0.005789"
$ws.Range("E16").Value = "'1.67%"
$ws.Range("D17").Value = "'3.466"
$ws.Range("E17").Value = "'0.32%"
$ws.Range("E18").Value = "'6.56%"
$ws.Range("E19").Value = "'1.15%"
$ws.Range("D20").Value = "'0.1330"
$ws.Range("E20").Value = "'-0.93%"
$ws.Range("D21").Value = "'5.105"
$ws.Range("E21").Value = "'-1.06%"
$ws.Range("D22").Value = "'0.2213"
$ws.Range("E22").Value = "'0.07%"
$ws.Range("D23").Value = "'0.04547"
$ws.Range("E23").Value = "'1.14%"
$ws.Range("D24").Value = "'0.001249"
$ws.Range("E24").Value = "'1.08%"
$ws.Range("D25").Value = "'0.004689"
$ws.Range("E25").Value = "'-3.50%"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("E26").Value = "'-3.78%"
$ws.Range("D27").Value = "'0.0004504"
$ws.Range("E27").Value = "'-5.33%"
$ws.Range("D39").Value = "'0.01962"
$ws.Range("E39").Value = "'7.39%"
$ws.Range("D40").Value = "'0.04889"
$ws.Range("E40").Value = "'4.29%"
$ws.Range("D41").Value = "'0.007601"
$ws.Range("E41").Value = "'-1.77%"
$ws.Range("D42").Value = "'0.1386"
$ws.Range("E42").Value = "'-0.03%"
$ws.Range("D43").Value = "'0.007799"
$ws.Range("E43").Value = "'0.50%"
$ws.Range("D44").Value = "'0.002108"
$ws.Range("E44").Value = "'-3.89%"
$ws.Range("D45").Value = "'0.01175"
$ws.Range("E45").Value = "'5.38%"
$ws.Range("D46").Value = "'0.00006744"
$ws.Range("E46").Value = "'7.40%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D48").Value = "'154.14"
$ws.Range("E48").Value = "'203.79%"
$ws.Range("D49").Value = "'0.001702"
$ws.Range("E49").Value = "'-10.60%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.07%"
